# UC005 - Listar Empenhos Pendentes
# From v1.0.3 to v1.1
#
# The "Clica para atribuir/desatribuir o registro a si mesmo" test step
# (previously under TC4, row 33) is swapped with the "Clica para realizar
# o empenho de uma diaria" / "Apresenta a tela de Registrar Empenho" test
# step (previously under TC3, row 26). The TC3/TC4 labels themselves stay
# where they are; only the step/expected-result content moves between the
# two test case blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Content currently living in the TC3 block (row 26)
$tc3Step = $ws.Range("B26").Value2
$tc3Expected = $ws.Range("D26").Value2

# Content currently living in the TC4 block (row 33)
$tc4Step = $ws.Range("B33").Value2
$tc4Expected = $ws.Range("D33").Value2

# Swap: TC3 block now gets the old TC4 content, and vice versa
$ws.Range("B26").Value2 = $tc4Step
$ws.Range("D26").Value2 = $tc4Expected

$ws.Range("B33").Value2 = $tc3Step
$ws.Range("D33").Value2 = $tc3Expected
